$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 12182.777
$ws.Range("I28").Value = 15177
$ws.Range("K28").Value = 15177
$ws.Range("M28").Value = -14692

$ws.Range("H69").Value = 62506410
$ws.Range("I69").Value = 5351.5
$ws.Range("K69").Value = 16054.5
$ws.Range("M69").Value = -15180.5

$ws.Range("H72").Value = 62506410
$ws.Range("I72").Value = 5351.5
$ws.Range("K72").Value = 48163.5
$ws.Range("M72").Value = -43795.5

$ws.Range("H116").Value = 4295.8667
$ws.Range("I116").Value = 3831.111
$ws.Range("K116").Value = 3831.111
$ws.Range("M116").Value = -389.1109999999999

$ws.Range("H137").Value = 6878.7915
$ws.Range("I137").Value = 8560.888999999999
$ws.Range("K137").Value = 25682.667
$ws.Range("M137").Value = -23132.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2427
$ws.Range("I2").Value = 841.7143
$ws.Range("K2").Value = 841.7143
$ws.Range("M2").Value = -728.7143

$ws.Range("H11").Value = 4102600.5
$ws.Range("J11").Value = 5000
$ws.Range("L11").Value = 5000
$ws.Range("N11").Value = -5288

$ws.Range("H32").Value = 1825.7261
$ws.Range("I32").Value = 1746.6285
$ws.Range("K32").Value = 1746.6285
$ws.Range("M32").Value = -1459.6285

$ws.Range("H45").Value = 1162.8
$ws.Range("I45").Value = 966.6667
$ws.Range("K45").Value = 966.6667
$ws.Range("M45").Value = -589.6667

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents() | Out-Null

$ws.Range("H61").Value = 3730.2896
$ws.Range("I61").Value = 3672.6453
$ws.Range("J61").Value = 3985.5715
$ws.Range("K61").Value = 3672.6453
$ws.Range("L61").Value = 3985.5715
$ws.Range("M61").Value = -3460.6453
$ws.Range("N61").Value = -4409.5715

$ws.Range("H74").Value = 3414
$ws.Range("I74").Value = 3437.9333
$ws.Range("J74").Value = 3294.3333
$ws.Range("K74").Value = 3437.9333
$ws.Range("L74").Value = 3294.3333
$ws.Range("M74").Value = -2563.9333
$ws.Range("N74").Value = -5042.3333

$ws.Range("H77").Value = 3414
$ws.Range("I77").Value = 3437.9333
$ws.Range("J77").Value = 3294.3333
$ws.Range("K77").Value = 17189.6665
$ws.Range("L77").Value = 16471.6665
$ws.Range("M77").Value = -12821.6665
$ws.Range("N77").Value = -25207.6665

$ws.Range("H116").Value = 2427
$ws.Range("I116").Value = 841.7143
$ws.Range("K116").Value = 841.7143
$ws.Range("M116").Value = 1452.2857

$ws.Range("H122").Value = 13891608
$ws.Range("I122").Value = 13891608
$ws.Range("K122").Value = 41674824
$ws.Range("M122").Value = -41672374

$ws.Range("H132").Value = 3355.75
$ws.Range("I132").Value = 3306.0435
$ws.Range("K132").Value = 9918.130500000001
$ws.Range("M132").Value = -7388.130500000001

$ws.Range("H136").Value = 3730.2896
$ws.Range("I136").Value = 3672.6453
$ws.Range("J136").Value = 3985.5715
$ws.Range("K136").Value = 11017.9359
$ws.Range("L136").Value = 11956.7145
$ws.Range("M136").Value = -8467.9359
$ws.Range("N136").Value = -17056.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2427
$ws.Range("I3").Value = 841.7143
$ws.Range("K3").Value = 841.7143
$ws.Range("M3").Value = -727.7143

$ws.Range("H134").Value = 2541.0476
$ws.Range("I134").Value = 2555.2778
$ws.Range("K134").Value = 7665.8334
$ws.Range("M134").Value = -5130.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1470.0869
$ws.Range("I16").Value = 1350.7222
$ws.Range("K16").Value = 1350.7222
$ws.Range("M16").Value = -1063.7222

$ws.Range("H22").Value = 1122.7
$ws.Range("I22").Value = 849.1429000000001
$ws.Range("J22").Value = 1761
$ws.Range("K22").Value = 849.1429000000001
$ws.Range("L22").Value = 1761
$ws.Range("M22").Value = -499.1429000000001
$ws.Range("N22").Value = -2461

$ws.Range("H31").Value = 3536.8518
$ws.Range("I31").Value = 1603.8422
$ws.Range("J31").Value = 8127.75
$ws.Range("K31").Value = 1603.8422
$ws.Range("L31").Value = 8127.75
$ws.Range("M31").Value = -1308.8422
$ws.Range("N31").Value = -8717.75

$ws.Range("H34").Value = 3536.8518
$ws.Range("I34").Value = 1603.8422
$ws.Range("J34").Value = 8127.75
$ws.Range("K34").Value = 1603.8422
$ws.Range("L34").Value = 8127.75
$ws.Range("M34").Value = -1401.8422
$ws.Range("N34").Value = -8531.75

$ws.Range("H51").Value = 20022.5

$ws.Range("H61").Value = 20022.5

$ws.Range("H113").Value = 1470.0869
$ws.Range("I113").Value = 1350.7222
$ws.Range("K113").Value = 1350.7222
$ws.Range("M113").Value = 819.2778000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 245
$ws.Range("I34").Value = 245
$ws.Range("K34").Value = 735
$ws.Range("M34").Value = -651

$ws.Range("H131").Value = 1627767.2
$ws.Range("J131").Value = 2780285.5
$ws.Range("L131").Value = 8340856.5
$ws.Range("N131").Value = -8350936.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 1346499.9
$ws.Range("J21").Value = 2673000
$ws.Range("L21").Value = 2673000
$ws.Range("N21").Value = -2673346

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents() | Out-Null

$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents() | Out-Null

$ws.Range("H30").Value = 1346499.9
$ws.Range("J30").Value = 2673000
$ws.Range("L30").Value = 2673000
$ws.Range("N30").Value = -2673210

$ws.Range("H126").Value = 6889.4287
$ws.Range("I126").Value = 8476.333000000001
$ws.Range("J126").Value = 5699.25
$ws.Range("K126").Value = 25428.999
$ws.Range("L126").Value = 17097.75
$ws.Range("M126").Value = -22958.999
$ws.Range("N126").Value = -22037.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4356.8
$ws.Range("I136").Value = 4154.8623
$ws.Range("K136").Value = 12464.5869
$ws.Range("M136").Value = -9914.586899999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11980
$ws.Range("I62").Value = 11900
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 11900
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -11276
$ws.Range("N62").Value = -13248

$ws.Range("H65").Value = 11980
$ws.Range("I65").Value = 11900
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 59500
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -56380
$ws.Range("N65").Value = -66240

$ws.Range("H113").Value = 946.4054
$ws.Range("I113").Value = 833.75
$ws.Range("K113").Value = 2501.25
$ws.Range("M113").Value = -331.25

$ws.Range("H126").Value = 3178.3333
$ws.Range("I126").Value = 3191.0715
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 9573.2145
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -7103.2145
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 5829.6665
$ws.Range("I132").Value = 1907.05
$ws.Range("K132").Value = 5721.15
$ws.Range("M132").Value = -3191.15

$ws.Range("H136").Value = 2708
$ws.Range("I136").Value = 2291.2727
$ws.Range("K136").Value = 6873.8181
$ws.Range("M136").Value = -6873.8181
